$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in attendance marks for row 27 (Хасбулатов Магомед): columns C-F get 5
$ws.Range("C27").Value = 5
$ws.Range("D27").Value = 5
$ws.Range("E27").Value = 5
$ws.Range("F27").Value = 5

# Update the frozen-pane view so the visible top-left cell and the active
# selection reflect scrolling down to row 27 (matching the diff's sheetView
# changes: topLeftCell C4 -> C8, active selection G6 -> G27)
$ws.Range("G27").Select()
$activeWindow = $excel.ActiveWindow
$activeWindow.ScrollRow = 8
